# Fix the Ra_Stock_5 radium concentration on the Parameters sheet.
$wb = $excel.ActiveWorkbook

$paramWs = $wb.Worksheets.Item("Parameters")
$paramWs.Range("B6").Value = 1407
$paramWs.Range("C6").Value = 62
$paramWs.Range("C6").NumberFormat = "0.00E+00"

# Update the stale selections left over on other sheets.
$bottleWs = $wb.Worksheets.Item("Bottle Results")
$bottleWs.Range("S19").Select()

$avgWs = $wb.Worksheets.Item("Averaged Results")
$avgWs.Range("B7").Select()

# Leave "Parameters" as the active/selected sheet & cell.
$paramWs.Activate()
$paramWs.Range("B6:C6").Select()
